$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Note" text (L column) for existing data rows 6..35 and
#     update the G column (Value) figures with the recalculated series ---

$newNote = "Sum of Mpkm for all modes. Estimated values: Trams/Buses/Trolley estimated for <1998, Soft mobility <1994, public boats<2007."

$gValues = @{
    6  = 95133.808337202412
    7  = 96231.008337202395
    8  = 95600.6083372024
    9  = 95686.508337202409
    10 = 96024.337490088757
    11 = 95263.511635763178
    12 = 96729.547851711759
    13 = 97678.740304778476
    14 = 98906.521400979342
    15 = 100378.24041044943
    16 = 102218.27989685429
    17 = 103621.75646017618
    18 = 105514.1177275793
    19 = 106806.65015753482
    20 = 108114.07200453606
    21 = 109778.24511842233
    22 = 111219.815829339
    23 = 113202.1099954244
    24 = 115644.9614315371
    25 = 119244.17143367991
    26 = 122175.77117423796
    27 = 123446.93057644294
    28 = 124817.69570135229
    29 = 126348.64930924743
    30 = 128261.34407471852
    31 = 130084.58177608065
    32 = 132718.78289768533
    33 = 134748.7651726871
    34 = 135864.02593483773
    35 = 138324.18024710321
}

for ($r = 6; $r -le 35; $r++) {
    $ws.Cells.Item($r, 7).Value = $gValues[$r]
    $ws.Cells.Item($r, 12).Value = $newNote
}

# --- Append two new data rows for years 2020 and 2021 ---
# Columns: A Country, B Entity, C Type, D Parameter, E Year, F Flow,
#          G Value, H Unit, I Delete, J Reference, K Link, L Note

$ws.Cells.Item(36, 1).Value = "CHE"
$ws.Cells.Item(36, 2).Value = "passenger"
$ws.Cells.Item(36, 3).Value = "actual_flow"
$ws.Cells.Item(36, 4).Value = "annual"
$ws.Cells.Item(36, 5).Value = 2020
$ws.Cells.Item(36, 7).Value = 117969.10022506714
$ws.Cells.Item(36, 8).Value = "Mpkm"
$ws.Cells.Item(36, 10).Value = "OFS: Prestations du transport de personnes (PV-L), Statistique des transports publics (TP)"
$ws.Cells.Item(36, 11).Value = "https://www.bfs.admin.ch/bfs/fr/home/statistiques/mobilite-transports/enquetes/pv-l.html"
$ws.Cells.Item(36, 12).Value = $newNote

$ws.Cells.Item(37, 1).Value = "CHE"
$ws.Cells.Item(37, 2).Value = "passenger"
$ws.Cells.Item(37, 3).Value = "actual_flow"
$ws.Cells.Item(37, 4).Value = "annual"
$ws.Cells.Item(37, 5).Value = 2021
$ws.Cells.Item(37, 7).Value = 124656.13952624553
$ws.Cells.Item(37, 8).Value = "Mpkm"
$ws.Cells.Item(37, 10).Value = "OFS: Prestations du transport de personnes (PV-L), Statistique des transports publics (TP)"
$ws.Cells.Item(37, 11).Value = "https://www.bfs.admin.ch/bfs/fr/home/statistiques/mobilite-transports/enquetes/pv-l.html"
$ws.Cells.Item(37, 12).Value = $newNote

# Copy the "Link" column (K) cell format (small blue-ish hyperlink font style)
# from an existing row down onto the two new rows, then restore the values
# (PasteSpecial formats-only keeps the text intact, this is just to be safe).
$ws.Range("K6").Copy()
$ws.Range("K36:K37").PasteSpecial(-4122)
$ws.Cells.Item(36, 11).Value = "https://www.bfs.admin.ch/bfs/fr/home/statistiques/mobilite-transports/enquetes/pv-l.html"
$ws.Cells.Item(37, 11).Value = "https://www.bfs.admin.ch/bfs/fr/home/statistiques/mobilite-transports/enquetes/pv-l.html"

# --- Update the selection to reflect the new active range ---
$null = $ws.Range("L6:L37").Select()
